$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the reservation rows that had the "remise applied twice" mixup ---
# Row 4 user was "monji" -> should be "ahmed" (date/status unchanged)
$ws.Cells.Item(4, 1).Value = "ahmed"
# Row 5 user was "ahmed" -> should be "karoui" (date/status unchanged)
$ws.Cells.Item(5, 1).Value = "karoui"

# --- Append the new confirmed reservations (rows 6-13), all for "karoui" ---
$timestamps = @(
    "2025-03-06T12:34:49.917466400",
    "2025-03-06T12:40:07.667345",
    "2025-03-06T12:44:21.003296400",
    "2025-03-06T12:44:34.228014700",
    "2025-03-06T12:53:51.661579400",
    "2025-03-06T12:58:30.918849200",
    "2025-03-06T13:01:30.447345700",
    "2025-03-06T13:02:13.757866900"
)

$row = 6
foreach ($ts in $timestamps) {
    $ws.Cells.Item($row, 1).Value = "karoui"
    $ws.Cells.Item($row, 2).Value = $ts
    $ws.Cells.Item($row, 2).HorizontalAlignment = -4108
    $ws.Cells.Item($row, 3).Value = "Confirmé"
    # RGB CCFFCC (the same light-green already used by the workbook's
    # "fillId 7" / indexed-color-42 fill) packed as an OLE BGR integer.
    $ws.Cells.Item($row, 3).Interior.Color = 13434828
    $row = $row + 1
}

# Column B needs to widen to fit the long timestamp strings (bestFit)
$ws.Columns.Item(2).ColumnWidth = 30

Write-Host "done"
